$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '88.458.84'
Set-TextValue $ws 'E2' '  -2.35%  '

Set-TextValue $ws 'D3' '3.114.74'
Set-TextValue $ws 'E3' '  -1.89%  '

Set-TextValue $ws 'D4' '1.00'
Set-TextValue $ws 'E4' '  +0.40%  '

Set-TextValue $ws 'D5' '213.15'
Set-TextValue $ws 'E5' '  +0.34%  '

Set-TextValue $ws 'D6' '633.51'
Set-TextValue $ws 'E6' '  +3.01%  '

Set-TextValue $ws 'D7' '0.389'
Set-TextValue $ws 'E7' '  +0.07%  '

Set-TextValue $ws 'D8' '0.776'
Set-TextValue $ws 'E8' '  +13.20%  '

Set-TextValue $ws 'D9' '1.00'
Set-TextValue $ws 'E9' '  +0.25%  '

Set-TextValue $ws 'D10' '3.112.07'
Set-TextValue $ws 'E10' '  -1.73%  '

Set-TextValue $ws 'D11' '0.559'
Set-TextValue $ws 'E11' '  -2.46%  '

Set-TextValue $ws 'D12' '0.178'
Set-TextValue $ws 'E12' '  +1.47%  '

Set-TextValue $ws 'D13' '0.0000249'
Set-TextValue $ws 'E13' '  -0.82%  '

Set-TextValue $ws 'D14' '5.34'
Set-TextValue $ws 'E14' '  +2.61%  '

Set-TextValue $ws 'D15' '88.468.65'
Set-TextValue $ws 'E15' '  -2.03%  '

Set-TextValue $ws 'D16' '3.692.28'
Set-TextValue $ws 'E16' '  -1.48%  '

Set-TextValue $ws 'D17' '31.98'
Set-TextValue $ws 'E17' '  -2.06%  '

Set-TextValue $ws 'D18' '3.121.02'
Set-TextValue $ws 'E18' '  -1.80%  '

Set-TextValue $ws 'D19' '3.41'
Set-TextValue $ws 'E19' '  +4.42%  '

Set-TextValue $ws 'D20' '0.0000221'
Set-TextValue $ws 'E20' '  +18.50%  '

Set-TextValue $ws 'D21' '13.13'
Set-TextValue $ws 'E21' '  -1.87%  '

Set-TextValue $ws 'D22' '420.12'
Set-TextValue $ws 'E22' '  -3.21%  '

Set-TextValue $ws 'D23' '8.39'
Set-TextValue $ws 'E23' '  -1.49%  '

Set-TextValue $ws 'D24' '4.88'
Set-TextValue $ws 'E24' '  -4.02%  '

Set-TextValue $ws 'D25' '5.37'
Set-TextValue $ws 'E25' '  +5.24%  '

Set-TextValue $ws 'D26' '81.92'
Set-TextValue $ws 'E26' '  +10.28%  '

Set-TextValue $ws 'D27' '11.39'
Set-TextValue $ws 'E27' '  -3.21%  '

Set-TextValue $ws 'D28' '3.297.25'
Set-TextValue $ws 'E28' '  -2.00%  '

Set-TextValue $ws 'E29' '  +0.02%  '

Set-TextValue $ws 'E30' '  -0.13%  '

Set-TextValue $ws 'D31' '0.156'
Set-TextValue $ws 'E31' '  -7.76%  '

Set-TextValue $ws 'D32' '3.98'
Set-TextValue $ws 'E32' '  -5.74%  '

Set-TextValue $ws 'D33' '8.12'
Set-TextValue $ws 'E33' '  -3.98%  '

Set-TextValue $ws 'D34' '500.50'
Set-TextValue $ws 'E34' '  -5.79%  '

Set-TextValue $ws 'D35' '0.147'
Set-TextValue $ws 'E35' '  +15.96%  '

Set-TextValue $ws 'D36' '6.88'
Set-TextValue $ws 'E36' '  -0.50%  '

Set-TextValue $ws 'D37' '1.27'
Set-TextValue $ws 'E37' '  +1.97%  '

Set-TextValue $ws 'D38' '1.83'
Set-TextValue $ws 'E38' '  -1.97%  '

Set-TextValue $ws 'D39' '22.19'
Set-TextValue $ws 'E39' '  +1.61%  '

Set-TextValue $ws 'D40' '22.21'
Set-TextValue $ws 'E40' '  -0.25%  '

Set-TextValue $ws 'E41' '  +0.50%  '

Set-TextValue $ws 'E42' '  +0.12%  '

Set-TextValue $ws 'D43' '0.363'
Set-TextValue $ws 'E43' '  -3.08%  '

Set-TextValue $ws 'D44' '1.84'
Set-TextValue $ws 'E44' '  -3.52%  '

Set-TextValue $ws 'B45' 'Monero'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws 'D45' '145.86'
Set-TextValue $ws 'E45' '  -0.55%  '

Set-TextValue $ws 'B46' 'Stellar'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws 'D46' '0.131'
Set-TextValue $ws 'E46' '  +6.82%  '

Set-TextValue $ws 'D47' '43.60'
Set-TextValue $ws 'E47' '  -2.02%  '

Set-TextValue $ws 'B48' 'Hedera'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws 'D48' '0.0660'
Set-TextValue $ws 'E48' '  +12.54%  '

Set-TextValue $ws 'B49' 'Aave'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D49' '162.36'
Set-TextValue $ws 'E49' '  -5.90%  '

Set-TextValue $ws 'D50' '0.715'
Set-TextValue $ws 'E50' '  +2.18%  '

Set-TextValue $ws 'D51' '1.18'
Set-TextValue $ws 'E51' '  -3.96%  '
